# fix(contracts): ajustar bloque de firmas en plantilla propietario
#
# 1) For every body paragraph whose paragraph properties are exactly
#    spacing-after=120 twips (6pt) + justify ("both"), explicitly stamp
#    KeepWithNext=False and KeepTogether=False (these serialize to
#    <w:keepNext w:val=".../> <w:keepLines .../> in the saved OOXML).
# 2) Rework the signature block at the end of the document (the two
#    "_____" signature-line groups) so each line keeps with the next
#    line/page and uses the tighter spacing from the new template.

$d = $word.ActiveDocument

# wdParagraphAlignment
$wdAlignParagraphLeft = 0
$wdAlignParagraphJustify = 3

# --- 1) stamp keepNext/keepLines = False on the "spacing after=120 / jc=both" paragraphs ---
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $fmt = $p.Format
    if ($fmt.SpaceBefore -eq 0 -and $fmt.SpaceAfter -eq 6 -and $fmt.Alignment -eq $wdAlignParagraphJustify) {
        $fmt.KeepWithNext = $false
        $fmt.KeepTogether = $false
    }
}

# --- 2) signature block rework ---
# Paragraph indices (fixed - no insert/delete of paragraphs happens above):
#  87 "________________________"                       (Arrendador's line)
#  88 "[[ARRENDADOR.NOMBRE]]"
#  89 "Rut: [[ARRENDADOR.RUT]]"
#  90 "Arrendador"
#  91 "" (blank spacer, untouched)
#  92 "________________________"                       (Arrendataria's line)
#  93 "pp. [[ARRENDATARIA.REPRESENTANTE.NOMBRE]]"
#  94 "Rut: [[ARRENDATARIA.REPRESENTANTE.RUT]]"
#  95 "Arrendataria"

$pLine1 = $d.Paragraphs.Item(87)
$pLine1.Format.KeepWithNext = $true
$pLine1.Format.KeepTogether = $true
$pLine1.Format.SpaceBefore = 6
$pLine1.Format.SpaceAfter = 2
$pLine1.Format.Alignment = $wdAlignParagraphLeft

$pName1 = $d.Paragraphs.Item(88)
$pName1.Format.KeepWithNext = $true
$pName1.Format.KeepTogether = $true
$pName1.Format.SpaceAfter = 2.75
$pName1.Format.Alignment = $wdAlignParagraphJustify
$rName1 = $pName1.Range
$rName1.Font.Bold = $false
$rName1.Font.Size = 11

$pRut1 = $d.Paragraphs.Item(89)
$pRut1.Format.KeepWithNext = $true
$pRut1.Format.KeepTogether = $true
$pRut1.Format.SpaceAfter = 2.75
$pRut1.Format.Alignment = $wdAlignParagraphLeft

$pRole1 = $d.Paragraphs.Item(90)
$pRole1.Format.KeepWithNext = $true
$pRole1.Format.KeepTogether = $true
$pRole1.Format.SpaceAfter = 2.75
$pRole1.Format.Alignment = $wdAlignParagraphLeft

$pLine2 = $d.Paragraphs.Item(92)
$pLine2.Format.KeepWithNext = $true
$pLine2.Format.KeepTogether = $true
$pLine2.Format.SpaceBefore = 6
$pLine2.Format.SpaceAfter = 2
$pLine2.Format.Alignment = $wdAlignParagraphLeft

$pName2 = $d.Paragraphs.Item(93)
$pName2.Format.KeepWithNext = $true
$pName2.Format.KeepTogether = $true
$pName2.Format.SpaceAfter = 2.75
$pName2.Format.Alignment = $wdAlignParagraphLeft

$pRut2 = $d.Paragraphs.Item(94)
$pRut2.Format.KeepWithNext = $true
$pRut2.Format.KeepTogether = $true
$pRut2.Format.SpaceAfter = 2.75
$pRut2.Format.Alignment = $wdAlignParagraphLeft

$pRole2 = $d.Paragraphs.Item(95)
$pRole2.Format.KeepWithNext = $false
$pRole2.Format.KeepTogether = $true
$pRole2.Format.SpaceAfter = 2.75
$pRole2.Format.Alignment = $wdAlignParagraphLeft

Write-Output "done"
